$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.158448219299316
$ws.Range("B1").Value = 2.691727876663208
$ws.Range("C1").Value = 2.861679077148438
$ws.Range("D1").Value = 3.475243330001831
$ws.Range("E1").Value = 1.81864595413208
